# "updated codebook based on meeting"
# Content changes happen on the "time-based helpful" sheet (Duration / Starting rows);
# the other two sheets only get their shared-string indices shifted as a side effect
# (handled automatically by the engine), plus all three sheets' view/selection state
# moves to reflect where the author was working when they saved.

$wb = $excel.ActiveWorkbook

$wsBody = $wb.Worksheets.Item("body-based helpful")
$wsTime = $wb.Worksheets.Item("time-based helpful")
$wsOther = $wb.Worksheets.Item("otherwise")

# --- Content edits on "time-based helpful" ---
# Row 5 "Duration": clarify the definition now also covers breaks, and broaden the example.
$wsTime.Range("B5").Value = "Telling a person how much time is left in an exercise or a break"
$wsTime.Range("C5").Value = "so we have 10 seconds left, 5 more seconds; getting ready to go"

# Row 3 "Starting an exercise": replace the example text.
$wsTime.Range("C3").Value = "starting, let's go"

# --- View/selection state to match the saved workbook ---
$wsBody.Range("B6").Select()
$wsOther.Range("A8").Select()

$wsTime.Activate()
$wsTime.Range("C4").Select()
